# Fruta / hortaliza, semanal
# A new weekly price observation (fecha serial 44722) is inserted ahead of
# the existing rows, pushing the prior rows 14-17 down to 15-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14, shifting the existing rows 14-17 down to 15-18.
$ws.Rows(14).Insert()

# Populate the newly inserted row 14 with the new weekly data point.
$ws.Range("A14").Value = 6
$ws.Range("B14").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C14").Value = 'Metropolitana'
$ws.Range("D14").Value = 44722
$ws.Range("E14").Value = 13
$ws.Range("F14").Value = 100112035
$ws.Range("G14").Value = 'Bruselas (repollito)'
$ws.Range("H14").Value = 'Sin especificar'
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 150
$ws.Range("K14").Value = 18000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 18933
$ws.Range("N14").Value = '$/malla 15 kilos'
$ws.Range("O14").Value = 'Provincia de Quillota'
$ws.Range("P14").Value = 1262
$ws.Range("Q14").Value = 15
$ws.Range("R14").Value = 'Hortaliza'
